# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# The old scraper only pulled team statistics, not the season record, so
# this backfills the record for every player row using the team's actual
# 2018 season totals (89 wins, 73 losses, 0 ties).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, one column to the right of the existing "Unnamed: 28"
# column (A..AC already in use).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header formatting used by the rest of
# row 1 by copying the format from the adjacent existing header cell.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill in the season record for every data row (2-55) with the team's
# win/loss/tie totals for the season.
$lastRow = 55
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 89   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 73   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
